$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated results after CMC (QPA) run: columns G (Absorbed Intensity by QPA),
# H (# of photons), and I (Rate of exciton generation) for rows 2-8.

$ws.Range("G2").Value = 76633.0298535531
$ws.Range("H2").Value = [double]"2.120319961672114e+23"
$ws.Range("I2").Value = [double]"4.240639923344228e+21"

$ws.Range("G3").Value = 76778.36861840234
$ws.Range("H3").Value = [double]"2.124341265343697e+23"
$ws.Range("I3").Value = [double]"4.248682530687393e+21"

$ws.Range("G4").Value = 76935.71167704028
$ws.Range("H4").Value = [double]"2.128694709657436e+23"
$ws.Range("I4").Value = [double]"4.257389419314872e+21"

$ws.Range("G5").Value = 77104.97017277453
$ws.Range("H5").Value = [double]"2.133377836083127e+23"
$ws.Range("I5").Value = [double]"4.266755672166253e+21"

$ws.Range("G6").Value = 77286.04848674376
$ws.Range("H6").Value = [double]"2.138387998991583e+23"
$ws.Range("I6").Value = [double]"4.276775997983167e+21"

$ws.Range("G7").Value = 77478.84428596617
$ws.Range("H7").Value = [double]"2.143722366984066e+23"
$ws.Range("I7").Value = [double]"4.287444733968133e+21"

$ws.Range("G8").Value = 77683.24857529382
$ws.Range("H8").Value = [double]"2.149377924329784e+23"
$ws.Range("I8").Value = [double]"4.298755848659567e+21"
